{"js": "// Lattice-multiplication practice sheet: regenerate the 5x3 table of\n// exercises with a new set of two-digit multiplication problems.\n//\n// Each table cell holds a single paragraph / single run containing five\n// lines (separated by <w:br/>, i.e. a vertical-tab \"\\v\" in Office.js\n// range text):\n//   1) \"A x B\"                 the two factors\n//   2) \"  b1    b2\"            spaced-out digits of B\n//   3) \"  ----\"                separator\n//   4) \"a1|    |\"              tens digit of A\n//   5) \"a2|    |\"              units digit of A\n//\n// New factors (row-major, 5 rows x 3 columns):\nconst newCells = [\n  [\"89 x 45\", \"  4    5\", \"  ----\", \"8|    |\", \"9|    |\"],\n  [\"45 x 50\", \"  5    0\", \"  ----\", \"4|    |\", \"5|    |\"],\n  [\"93 x 96\", \"  9    6\", \"  ----\", \"9|    |\", \"3|    |\"],\n  [\"28 x 73\", \"  7    3\", \"  ----\", \"2|    |\", \"8|    |\"],\n  [\"69 x 13\", \"  1    3\", \"  ----\", \"6|    |\", \"9|    |\"],\n  [\"98 x 49\", \"  4    9\", \"  ----\", \"9|    |\", \"8|    |\"],\n  [\"29 x 83\", \"  8    3\", \"  ----\", \"2|    |\", \"9|    |\"],\n  [\"75 x 72\", \"  7    2\", \"  ----\", \"7|    |\", \"5|    |\"],\n  [\"91 x 47\", \"  4    7\", \"  ----\", \"9|    |\", \"1|    |\"],\n  [\"27 x 72\", \"  7    2\", \"  ----\", \"2|    |\", \"7|    |\"],\n  [\"51 x 24\", \"  2    4\", \"  ----\", \"5|    |\", \"1|    |\"],\n  [\"65 x 10\", \"  1    0\", \"  ----\", \"6|    |\", \"5|    |\"],\n  [\"96 x 49\", \"  4    9\", \"  ----\", \"9|    |\", \"6|    |\"],\n  [\"72 x 82\", \"  8    2\", \"  ----\", \"7|    |\", \"2|    |\"],\n  [\"96 x 19\", \"  1    9\", \"  ----\", \"9|    |\", \"6|    |\"],\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount,values\");\nawait context.sync();\n\nconst columns = 3;\nfor (let i = 0; i < newCells.length; i++) {\n  const row = Math.floor(i / columns);\n  const col = i % columns;\n  if (row >= table.rowCount) {\n    break;\n  }\n  const cell = table.getCell(row, col);\n  const paragraph = cell.body.paragraphs.getFirst();\n  const range = paragraph.getRange();\n  // Join with vertical-tab (\\v) -> becomes a <w:br/> line break, matching\n  // the original run layout; this replaces the run text in place, keeping\n  // the existing run formatting (e.g. sz=32) intact.\n  range.insertText(newCells[i].join(\"\\v\"), Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Lattice-multiplication practice sheet: regenerate the 5x3 table of\n# exercises with a new set of two-digit multiplication problems.\n#\n# Each table cell holds five lines separated by line breaks (a vertical\n# tab, `v, in a COM Range.Text assignment becomes a <w:br/>):\n#   1) \"A x B\"                 the two factors\n#   2) \"  b1    b2\"            spaced-out digits of B\n#   3) \"  ----\"                separator\n#   4) \"a1|    |\"              tens digit of A\n#   5) \"a2|    |\"              units digit of A\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$newCells = @(\n  \"89 x 45`v  4    5`v  ----`v8|    |`v9|    |\",\n  \"45 x 50`v  5    0`v  ----`v4|    |`v5|    |\",\n  \"93 x 96`v  9    6`v  ----`v9|    |`v3|    |\",\n  \"28 x 73`v  7    3`v  ----`v2|    |`v8|    |\",\n  \"69 x 13`v  1    3`v  ----`v6|    |`v9|    |\",\n  \"98 x 49`v  4    9`v  ----`v9|    |`v8|    |\",\n  \"29 x 83`v  8    3`v  ----`v2|    |`v9|    |\",\n  \"75 x 72`v  7    2`v  ----`v7|    |`v5|    |\",\n  \"91 x 47`v  4    7`v  ----`v9|    |`v1|    |\",\n  \"27 x 72`v  7    2`v  ----`v2|    |`v7|    |\",\n  \"51 x 24`v  2    4`v  ----`v5|    |`v1|    |\",\n  \"65 x 10`v  1    0`v  ----`v6|    |`v5|    |\",\n  \"96 x 49`v  4    9`v  ----`v9|    |`v6|    |\",\n  \"72 x 82`v  8    2`v  ----`v7|    |`v2|    |\",\n  \"96 x 19`v  1    9`v  ----`v9|    |`v6|    |\"\n)\n\n$columns = 3\n$rowCount = $t.Rows.Count\n\nfor ($i = 0; $i -lt $newCells.Count; $i++) {\n  $row = [math]::Floor($i / $columns) + 1\n  $col = ($i % $columns) + 1\n  if ($row -gt $rowCount) {\n    break\n  }\n  $cell = $t.Cell($row, $col)\n  # Assigning .Range.Text replaces the cell's run text in place while\n  # keeping the existing run formatting (e.g. sz=32) intact.\n  $cell.Range.Text = $newCells[$i]\n}\n\nWrite-Output \"done\"\n"}
